$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.554.52"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.762.83"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3837"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3408"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.05"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.139"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07396"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.49"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.346"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "1.763.46"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.039"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06662"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.34"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.374"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.77%  "
$ws.Range("D23").Value = "27.560.31"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.03"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.387"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.59"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.42%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.422"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.413"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "152.55"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.63"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "1.962.86"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.107"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.954"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.70"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02416"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6772"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.323"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06314"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2177"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.505"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -9.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.238"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.22"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6244"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.35%  "
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.78"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.077"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07371"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.146"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.61%  "
